# "moved metadata in ID Sawtooth data"
#
# The original workbook has a single sheet ("Sheet1") containing the ID
# Sawtooth data table (A1:H13 with headers) plus a stray two-cell
# "MetaData" / long-comment row tacked on at row 17 (A17/B17).
#
# This edit:
#   1. Renames "Sheet1" -> "ID data".
#   2. Removes the stray metadata row (row 17) from "ID data".
#   3. Adds a new "metadata" worksheet (placed right after "ID data") that
#      holds that same information as a proper little table with headers
#      ("dataset" / "comments") and a dataset id (55) alongside the long
#      comment text.
#   4. Leaves "ID data" as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the data sheet ------------------------------------------------
$idData = $wb.Worksheets.Item(1)
$idData.Name = "ID data"

# --- 2. Drop the stray metadata row out of the data table -------------------
$idData.Rows.Item(17).Delete()

# --- 3. Build the new "metadata" sheet --------------------------------------
$metadata = $wb.Worksheets.Add($null, $idData)
$metadata.Name = "metadata"

$metadata.Range("A1").Value = "dataset"
$metadata.Range("B1").Value = "comments"
$metadata.Range("A2").Value = 55
$metadata.Range("B2").Value = "Hatchery weir operated at Sawtooth Fish Hatchery. Weir is operated to collect spring/summer Chinook for the hatchery program. However, BT are caught at the facility. Operations have been consistent across the dataseries."
$metadata.Range("B2").Select() | Out-Null

# --- 4. Leave focus back on the data sheet -----------------------------------
$idData.Activate() | Out-Null
$idData.Range("A1").Select() | Out-Null
